$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "720×2=" "646×5="
Replace-Text "259×9=" "717×5="
Replace-Text "496×2=" "683×5="
Replace-Text "755×6=" "291×8="
Replace-Text "779×4=" "947×4="
Replace-Text "702×7=" "653×2="
Replace-Text "378×2=" "665×6="
Replace-Text "122×3=" "272×9="
Replace-Text "898×6=" "620×5="
Replace-Text "191×7=" "130×7="
Replace-Text "349×7=" "514×2="
Replace-Text "661×3=" "911×8="
Replace-Text "130×4=" "482×8="
Replace-Text "411×6=" "438×9="
Replace-Text "799×4=" "844×6="
Replace-Text "446×6=" "307×2="
Replace-Text "543×8=" "732×2="
Replace-Text "602×2=" "370×4="
Replace-Text "424×8=" "830×4="
Replace-Text "939×8=" "400×2="
Replace-Text "943×2=" "616×5="
Replace-Text "650×7=" "725×8="
Replace-Text "578×3=" "339×2="
Replace-Text "480×8=" "197×6="
Replace-Text "855×8=" "848×7="
